$d = $word.ActiveDocument

# --- Step 1: remove the standalone "Meta description" paragraph that
#     directly follows the title heading ("Meta description: Read our
#     review ..."). ---
$metaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Meta description*") {
        $metaIndex = $i
        break
    }
}
if ($metaIndex -gt 0) {
    $d.Paragraphs.Item($metaIndex).Range.Delete() | Out-Null
}

# --- Step 2: locate the closing paragraph that still holds the old
#     "Create a Feature Image Prompt..." copy, and insert a brand new
#     bold paragraph ("Play free El Dorado the City of Gold online slot
#     game") immediately before it. ---
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Create a Feature Image Prompt*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -gt 1) {
    $pPrev = $d.Paragraphs.Item($targetIndex - 1)
    # A position landing inside the existing text of the preceding
    # paragraph (rather than right on the paragraph boundary) makes
    # InsertXML insert a clean, independent sibling <w:p> right after
    # that paragraph instead of merging runs into it.
    $pos = $pPrev.Range.Start + 1
    $insertionPoint = $d.Range($pos, $pos)
    $newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play free El Dorado the City of Gold online slot game</w:t></w:r></w:p>'
    $insertionPoint.InsertXML($newParaXml) | Out-Null
}

# --- Step 3: swap out the old "Create a Feature Image Prompt..." body
#     copy for the new meta-description sentence (keeps the italic run
#     formatting that is already on that paragraph). ---
$oldText = 'Create a Feature Image Prompt: El Dorado the City of Gold Design an eye-catching feature image for "El Dorado The City of Gold" that would capture the attention of online slot game players. The image should accurately represent the game''s theme and graphics, with a playful and cartoonish style. The image should showcase a Maya warrior wearing glasses, who has discovered a hidden treasure of gold and jewels in the jungle. The use of bright colors, intricate details, and fun symbols such as toucans, monkeys, and tree frogs, should be incorporated in the design to create an exciting visual experience. Make sure the image highlights the mythical appeal of El Dorado, and should be able to grab the viewer''s attention at a glance.'
$newText = 'Read our review of El Dorado the City of Gold online slot game and play for free today.'
$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null
